$d = $word.ActiveDocument
$s = $d.Styles("Hyperlink")
$s.Font.ColorIndex = 0
$s.Font.Underline = 0
